# Refitting NCDEs to individual patients (for manuscript figure)
#
# Adds a new "Label" column (H) that marks each patient row as
# Control (0) or MDD (1), and refreshes a handful of recomputed
# D/E/F values coming from the refit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column H (copy the bold/bordered header formatting from G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Label"

# Updated (refit) numeric values in D/E/F for a few rows
$ws.Range("D3").Value = [double]"1.817676254734765E-16"
$ws.Range("E3").Value = [double]"1.817676254734765E-16"

$ws.Range("D5").Value = [double]"0.1427025992007107"
$ws.Range("E5").Value = [double]"0.1427025992007107"

$ws.Range("D10").Value = [double]"0.4356345699425944"
$ws.Range("E10").Value = [double]"0.5643654300574057"

$ws.Range("F11").Value = [double]"1632.057006835938"

# New "Label" column values: 0 = Control, 1 = MDD
$labels = @(0, 0, 0, 0, 0, 1, 1, 1, 1, 1, 0, 0, 0, 0, 0, 1, 1, 1, 1, 1)
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $labels[$i]
}
